$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C46").Value = 7734
$ws.Range("C47:C59").Value = 7343
$ws.Range("C60:C157").Value = 7293
